$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at F (pushes old F "formername" / G "formername2" -> G / H) ---
$ws.Columns("F:F").Insert()

# Re-apply the custom width (25.5) across B:F (was B:E before the new column existed)
$ws.Columns("B:F").ColumnWidth = 24.7

# --- Header row ---
$ws.Range("F1").Value = "altname5"

# --- New / updated data cells ---

# Row 49 (COD / DR Congo): new altname5 value, inserted before the shifted "Zaire"
$ws.Range("F49").Value = "Congo (Democratic Republic of the)"

# Row 140 (KOR / South Korea): new altname3 value
$ws.Range("D140").Value = "Korea (Republic of)"

# Row 141 (MDA / Moldova): new altname2 value
$ws.Range("C141").Value = "Moldova (Republic of)"

# Row 167 (SWZ / Swaziland): new altname value
$ws.Range("B167").Value = "Eswatini (Kingdom of)"

# Row 186 (TZA / Tanzania): new altname2 value
$ws.Range("C186").Value = "Tanzania (United Republic of)"

# Row 187 (USA): new altname2 value
$ws.Range("C187").Value = "United States"

# Row 197 (PSE / Palestine): new altname2 value
$ws.Range("C197").Value = "Palestine, State of"

# New row 199 (Hong Kong)
$ws.Range("A199").Value = "HKG"
$ws.Range("B199").Value = "Hong Kong, China (SAR)"

# --- View state: selection moves to F49 ---
$ws.Range("F49").Select()
